$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 4724.9487
$ws.Range("I100").Value = 3081.0715
$ws.Range("J100").Value = 5645.52
$ws.Range("K100").Value = 3081.0715
$ws.Range("L100").Value = 5645.52
$ws.Range("M100").Value = -2540.0715
$ws.Range("N100").Value = -6727.52

# Row 132
$ws.Range("H132").Value = 4206.0513
$ws.Range("I132").Value = 2045.1765
$ws.Range("K132").Value = 6135.529500000001
$ws.Range("M132").Value = -3605.529500000001

# Row 137
$ws.Range("H137").Value = 24566.8
$ws.Range("I137").Value = 46522.434
$ws.Range("J137").Value = 1613.1818
$ws.Range("K137").Value = 139567.302
$ws.Range("L137").Value = 4839.5454
$ws.Range("M137").Value = -137017.302
$ws.Range("N137").Value = -9939.545399999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 409.35715
$ws.Range("I2").Value = 410.07693
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 410.07693
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -297.07693
$ws.Range("N2").Value = -626

# Row 44
$ws.Range("H44").Value = 8124.8335
$ws.Range("I44").Value = 7500
$ws.Range("J44").Value = 8249.799999999999
$ws.Range("K44").Value = 7500
$ws.Range("L44").Value = 8249.799999999999
$ws.Range("M44").Value = -7012
$ws.Range("N44").Value = -9225.799999999999

# Row 55
$ws.Range("H55").Value = 13500
$ws.Range("J55").Value = 13500
$ws.Range("L55").Value = 13500
$ws.Range("N55").Value = -14130

# Row 63
$ws.Range("H63").Value = 4472.6
$ws.Range("I63").Value = 4799.0835
$ws.Range("J63").Value = 3166.6667
$ws.Range("K63").Value = 4799.0835
$ws.Range("L63").Value = 3166.6667
$ws.Range("M63").Value = -4113.0835
$ws.Range("N63").Value = -4538.6667

# Row 66
$ws.Range("H66").Value = 4472.6
$ws.Range("I66").Value = 4799.0835
$ws.Range("J66").Value = 3166.6667
$ws.Range("K66").Value = 23995.4175
$ws.Range("L66").Value = 15833.3335
$ws.Range("M66").Value = -20563.4175
$ws.Range("N66").Value = -22697.3335

# Row 116
$ws.Range("H116").Value = 409.35715
$ws.Range("I116").Value = 410.07693
$ws.Range("J116").Value = 400
$ws.Range("K116").Value = 410.07693
$ws.Range("L116").Value = 400
$ws.Range("M116").Value = 1883.92307
$ws.Range("N116").Value = -4988

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 409.35715
$ws.Range("I3").Value = 410.07693
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 410.07693
$ws.Range("L3").Value = 400
$ws.Range("M3").Value = -296.07693
$ws.Range("N3").Value = -628

# Row 20
$ws.Range("H20").Value = 1834.4546
$ws.Range("I20").Value = 1862.6
$ws.Range("J20").Value = 1774.1428
$ws.Range("K20").Value = 1862.6
$ws.Range("L20").Value = 1774.1428
$ws.Range("M20").Value = -1615.6
$ws.Range("N20").Value = -2268.1428

# Row 86
$ws.Range("H86").Value = 5116.357
$ws.Range("I86").Value = 4991.8667
$ws.Range("J86").Value = 5260
$ws.Range("K86").Value = 4991.8667
$ws.Range("L86").Value = 5260
$ws.Range("M86").Value = -3868.8667
$ws.Range("N86").Value = -7506

# Row 89
$ws.Range("H89").Value = 5116.357
$ws.Range("I89").Value = 4991.8667
$ws.Range("J89").Value = 5260
$ws.Range("K89").Value = 24959.3335
$ws.Range("L89").Value = 26300
$ws.Range("M89").Value = -19343.3335
$ws.Range("N89").Value = -37532

# Row 94
$ws.Range("H94").Value = 840.56665
$ws.Range("I94").Value = 707.617
$ws.Range("K94").Value = 707.617
$ws.Range("M94").Value = -256.617

# Row 116
$ws.Range("H116").Value = 19950
$ws.Range("J116").Value = 19950
$ws.Range("L116").Value = 19950
$ws.Range("N116").Value = -29128

# Row 140
$ws.Range("H140").Value = 20776.25
$ws.Range("J140").Value = 20776.25
$ws.Range("L140").Value = 20776.25
$ws.Range("N140").Value = -31136.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3966.6667
$ws.Range("I16").Value = 5500
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 5500
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -5213
$ws.Range("N16").Value = -1474

# Row 86
$ws.Range("H86").Value = 3924.7144
$ws.Range("I86").Value = 2200.1667
$ws.Range("J86").Value = 7028.9
$ws.Range("K86").Value = 2200.1667
$ws.Range("L86").Value = 7028.9
$ws.Range("M86").Value = -1077.1667
$ws.Range("N86").Value = -9274.9

# Row 89
$ws.Range("H89").Value = 3924.7144
$ws.Range("I89").Value = 2200.1667
$ws.Range("J89").Value = 7028.9
$ws.Range("K89").Value = 11000.8335
$ws.Range("L89").Value = 35144.5
$ws.Range("M89").Value = -5384.833500000001
$ws.Range("N89").Value = -46376.5

# Row 113
$ws.Range("H113").Value = 3966.6667
$ws.Range("I113").Value = 5500
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 5500
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = -3330
$ws.Range("N113").Value = -5240

# Row 134
$ws.Range("H134").Value = 1819.9025
$ws.Range("I134").Value = 1071.9524
$ws.Range("J134").Value = 2605.25
$ws.Range("K134").Value = 3215.857199999999
$ws.Range("L134").Value = 7815.75
$ws.Range("M134").Value = -680.8571999999995
$ws.Range("N134").Value = -12885.75

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 40.4
$ws.Range("J12").Value = 45.53846
$ws.Range("L12").Value = 136.61538
$ws.Range("N12").Value = -482.61538

# Row 107
$ws.Range("H107").Value = 248.67567
$ws.Range("I107").Value = 176.72
$ws.Range("J107").Value = 398.58334
$ws.Range("K107").Value = 530.16
$ws.Range("L107").Value = 1195.75002
$ws.Range("M107").Value = 1389.84
$ws.Range("N107").Value = -5035.750019999999

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1291.45
$ws.Range("I122").Value = 1287.7858
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 3863.3574
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -1413.3574
$ws.Range("N122").Value = -8800

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2092.7334
$ws.Range("I7").Value = 1909.2106
$ws.Range("J7").Value = 2409.7273
$ws.Range("K7").Value = 1909.2106
$ws.Range("L7").Value = 2409.7273
$ws.Range("M7").Value = -1797.2106
$ws.Range("N7").Value = -2633.7273

# Row 126
$ws.Range("H126").Value = 2092.7334
$ws.Range("I126").Value = 1909.2106
$ws.Range("J126").Value = 2409.7273
$ws.Range("K126").Value = 5727.6318
$ws.Range("L126").Value = 7229.1819
$ws.Range("M126").Value = -3257.6318
$ws.Range("N126").Value = -12169.1819

# Row 132
$ws.Range("H132").Value = 5297.2983
$ws.Range("I132").Value = 1725.5428
$ws.Range("J132").Value = 10979.637
$ws.Range("K132").Value = 5176.6284
$ws.Range("L132").Value = 32938.911
$ws.Range("M132").Value = -2646.6284
$ws.Range("N132").Value = -37998.911

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 18183678
$ws.Range("I136").Value = 35715820
$ws.Range("J136").Value = 2196.2222
$ws.Range("K136").Value = 107147460
$ws.Range("L136").Value = 6588.6666
$ws.Range("M136").Value = -107144910
$ws.Range("N136").Value = -11688.6666
